$wb = $excel.ActiveWorkbook

# --- Sheet "DBD": add new "Updated By" note column. ---
$dbd = $wb.Worksheets.Item("DBD")

$dbd.Range("H1").Value = "更新BY"
$dbd.Range("H2").Value = "L7205-五類資產分類上傳轉檔作業"

# --- Sheet "SP": stored-procedure parameter list now also carries jobTxSeq. ---
$sp = $wb.Worksheets.Item("SP")
$sp.Range("B2").Value = "int tbsdyf,  String empNo,  String jobTxSeq"
$sp.Range("B2").Select()

# --- Back on "DBD": update column type labels DATE -> TIMESTAMP, and drop the
#     trailing blank rows. ---
$dbd.Range("D13").Value = "TIMESTAMP"
$dbd.Range("D15").Value = "TIMESTAMP"

$dbd.Range("A17:A19").EntireRow.Delete()

$dbd.Range("D15").Select()
